$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.6311703333333334
$ws.Range("H2").Value = 1.893511
$ws.Range("M2").Value = 0.3284223333333333
$ws.Range("N2").Value = 0.985267
$ws.Range("O2").Value = 0.04541528350839906
$ws.Range("P2").Value = 0.04541528350839906
$ws.Range("Q2").Value = 0.2072904336041111
$ws.Range("R2").Value = 1.865613902437
$ws.Range("S2").Value = 0.04541528350839906
$ws.Range("T2").Value = 0.04541528350839906

# Row 3
$ws.Range("G3").Value = 0.6311703333333334
$ws.Range("H3").Value = 1.893511
$ws.Range("M3").Value = 4.062688333333334
$ws.Range("O3").Value = 0.5618014481290817
$ws.Range("P3").Value = 0.5618014481290816
$ws.Range("Q3").Value = 2.564248349579445
$ws.Range("R3").Value = 23.078235146215
$ws.Range("S3").Value = 0.5618014481290817
$ws.Range("T3").Value = 0.5618014481290816

# Row 4
$ws.Range("G4").Value = 0.6311703333333334
$ws.Range("H4").Value = 1.893511
$ws.Range("M4").Value = 2.840427
$ws.Range("N4").Value = 8.521281
$ws.Range("O4").Value = 0.3927832683625193
$ws.Range("P4").Value = 0.3927832683625193
$ws.Range("Q4").Value = 1.792793256399
$ws.Range("R4").Value = 16.135139307591
$ws.Range("S4").Value = 0.3927832683625193
$ws.Range("T4").Value = 0.3927832683625193
